$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.012.92"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.29%  "
$ws.Range("D3").Value = "'2.102.70"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.04%  "
$ws.Range("E4").Value = "  -0.99%  "
$ws.Range("D5").Value = "'346.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.42%  "
$ws.Range("E6").Value = "  -0.98%  "
$ws.Range("D7").Value = "'0.5155"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.46%  "
$ws.Range("D8").Value = "'0.4432"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.22%  "
$ws.Range("D9").Value = "'0.09402"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.06%  "
$ws.Range("D10").Value = "'52.46"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.17%  "
$ws.Range("D11").Value = "'1.175"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.27%  "
$ws.Range("D13").Value = "'2.103.59"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.51%  "
$ws.Range("D14").Value = "'6.751"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.74%  "
$ws.Range("D15").Value = "'8.185"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.55%  "
$ws.Range("D16").Value = "'99.69"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.34%  "
$ws.Range("D17").Value = "'0.00001162"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.36%  "
$ws.Range("E18").Value = "  -1.06%  "
$ws.Range("D19").Value = "'20.70"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.59%  "
$ws.Range("D20").Value = "'0.06684"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.35%  "
$ws.Range("E21").Value = "  -0.96%  "
$ws.Range("D22").Value = "'6.237"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.72%  "
$ws.Range("D23").Value = "'30.093.28"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.29%  "
$ws.Range("D24").Value = "'12.67"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.72%  "
$ws.Range("D25").Value = "'2.328"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.37%  "
$ws.Range("D26").Value = "'2.344.15"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.60%  "
$ws.Range("D27").Value = "'22.07"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.04%  "
$ws.Range("D28").Value = "'2.557"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.26%  "
$ws.Range("D29").Value = "'162.87"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.59%  "
$ws.Range("D30").Value = "'133.55"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.73%  "
$ws.Range("D31").Value = "'1.171"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.59%  "
$ws.Range("E32").Value = "  -1.92%  "
$ws.Range("D33").Value = "'1.644"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.06%  "
$ws.Range("D34").Value = "'6.246"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.32%  "
$ws.Range("D35").Value = "'3.954"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.08%  "
$ws.Range("D36").Value = "'6.190"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.37%  "
$ws.Range("D37").Value = "'10.17"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.25%  "
$ws.Range("D38").Value = "'0.02572"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.39%  "
$ws.Range("D39").Value = "'0.06775"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.79%  "
$ws.Range("D40").Value = "'0.2287"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.09%  "
$ws.Range("D41").Value = "'12.57"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.82%  "
$ws.Range("D42").Value = "'0.6931"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("E43").Value = "  +3.37%  "
$ws.Range("D44").Value = "'0.6679"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.77%  "
$ws.Range("D45").Value = "'14.23"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.03%  "
$ws.Range("D46").Value = "'2.298"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.83%  "
$ws.Range("D47").Value = "'3.635"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.95%  "
$ws.Range("D48").Value = "'0.00000000352"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.29%  "
$ws.Range("E49").Value = "  -3.30%  "
$ws.Range("D50").Value = "'82.18"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.11%  "
$ws.Range("D51").Value = "'0.07205"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.46%  "
